$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update existing shared string text used by G25 (Purpose of wget...)
# ------------------------------------------------------------------
$ws.Range("G25").Value = 'Purpose of wget, where to run command, optimizations, line continuations'

# ------------------------------------------------------------------
# 2. Populate new data rows 26-41 (columns A,C,D,E,F,G)
#    (column B - the 'Link' hyperlink cells - handled in step 3)
# ------------------------------------------------------------------
# Row 26
$ws.Range("A26").Value = 25
$ws.Range("C26").Value = 'Antlr quickstart'
$ws.Range("D26").Value = 'Quick Start (Antlr)'
$ws.Range("E26").Value = 'Overview'
$ws.Range("F26").Value = 'Succinct'
$ws.Range("G26").Value = 'grun meaning, curl, export, executing Java path'

# Row 27
$ws.Range("A27").Value = 26
$ws.Range("C27").Value = 'Antlr quickstart'
$ws.Range("D27").Value = 'Samples (Antlr)'
$ws.Range("E27").Value = 'Demonstrate API'
$ws.Range("F27").Value = 'Succinct'
$ws.Range("G27").Value = '"grammar" purpose, RegEx, grun command, where to paste code'

# Row 28
$ws.Range("A28").Value = 27
$ws.Range("C28").Value = 'Antlr setup'
$ws.Range("D28").Value = 'Getting Started with ANTLR v4'
$ws.Range("E28").Value = 'Show Process'
$ws.Range("F28").Value = 'Verbose'
$ws.Range("G28").Value = 'Save to /usr/local/lib on OSX, UNIX = OSX, adding to CLASSPATH, -Xmx500M, r command line opt'

# Row 29
$ws.Range("A29").Value = 28
$ws.Range("C29").Value = 'IntelliJ setup'
$ws.Range("D29").Value = 'Configuring Module Dependencies and Libraries'
$ws.Range("E29").Value = 'Teach'
$ws.Range("F29").Value = 'Verbose'
$ws.Range("G29").Value = 'terms: module library, module dependencies, project structure'

# Row 30
$ws.Range("A30").Value = 29
$ws.Range("C30").Value = 'Antlr guide'
$ws.Range("D30").Value = 'Grammar Lexicon'
$ws.Range("E30").Value = 'Standards'
$ws.Range("F30").Value = 'Verbose'
$ws.Range("G30").Value = 'terms: actions, Javadocs, Unicode, escape, reserved words'

# Row 31
$ws.Range("A31").Value = 30
$ws.Range("C31").Value = 'Antlr guide'
$ws.Range("D31").Value = 'Grammar structure'
$ws.Range("E31").Value = 'Standards'
$ws.Range("F31").Value = 'Verbose'
$ws.Range("G31").Value = 'terms: channels, lexer, parser, actions, depth-first search; | (or)'

# Row 32
$ws.Range("A32").Value = 31
$ws.Range("C32").Value = 'nginx HTTPS'
$ws.Range("D32").Value = 'How To Create an SSL Certificate on Nginx for Ubuntu 14.04'
$ws.Range("E32").Value = 'Teach'
$ws.Range("F32").Value = 'Verbose'
$ws.Range("G32").Value = 'nginx ''listen 80'', enable 443 only, Beast Attack, certificate authority'

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("C33").Value = 'Java HTTP server'
$ws.Range("D33").Value = 'Hava a simple HTTP server'
$ws.Range("E33").Value = 'Teach'
$ws.Range("F33").Value = 'Succinct'
$ws.Range("G33").Value = 'concept: JDK, HTTP server; throws Exception, InetSocketAddress, createContext, HttpHandler,  HttpExchange methods?, setExecutor, display PDF, naming class file'

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("C34").Value = 'Java write JSON'
$ws.Range("D34").Value = 'JSON.Simple Example – Read And Write JSON'
$ws.Range("E34").Value = 'Teach'
$ws.Range("F34").Value = 'Succinct'
$ws.Range("G34").Value = 'terms: JSON; Maven setup, pom.xml, parse pure string how, where to insert pom dependency, what the rest says, empty JAR, including SRC, catch ParseException, local jars as deps'

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("C35").Value = 'Maven POM definition'
$ws.Range("D35").Value = 'Introduction to the POM'
$ws.Range("E35").Value = 'Teach'
$ws.Range("F35").Value = 'Verbose'
$ws.Range("G35").Value = 'maven, POM pieces: repository layout, url, updatePolicy, etc., FQAN'

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("C36").Value = 'Maven build project'
$ws.Range("D36").Value = 'Maven Tutorial'
$ws.Range("E36").Value = 'Gateway'
$ws.Range("F36").Value = 'Verbose'
$ws.Range("G36").Value = 'proxy access, Eclipse IDE, web application, Tomcat'

# Row 37
$ws.Range("A37").Value = 36
$ws.Range("C37").Value = 'Maven build project'
$ws.Range("D37").Value = 'How To Build Project With Maven'
$ws.Range("E37").Value = 'Teach (Fragment)'
$ws.Range("F37").Value = 'Succinct'
$ws.Range("G37").Value = 'maven, source location, install maven, war, packaging options'

# Row 38
$ws.Range("A38").Value = 37
$ws.Range("C38").Value = 'Maven build project'
$ws.Range("D38").Value = 'How To Create A Java Project With Maven'
$ws.Range("E38").Value = 'Teach'
$ws.Range("F38").Value = 'Verbose'
$ws.Range("G38").Value = '{} means custom param, /src/ vs. src/, business logic, -cp, Eclipse IDE, example archetypes, non-empty dir init'

# Row 39
$ws.Range("A39").Value = 38
$ws.Range("C39").Value = 'Maven run class'
$ws.Range("D39").Value = '3 ways to run Java main from Maven'
$ws.Range("E39").Value = 'Teach'
$ws.Range("F39").Value = 'Verbose'
$ws.Range("G39").Value = 'mvn installed, open Terminal, exec:java goal, maven phases, maven profiles, exec:java missing'

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("C40").Value = 'Maven local JAR'
$ws.Range("D40").Value = 'How to add local jar files in maven project?'
$ws.Range("F40").Value = 'Succinct'
$ws.Range("G40").Value = 'Don''t copy-paste all, example group-id, example artifact, path cannot include ''~'', '

# Row 41
$ws.Range("A41").Value = 40
$ws.Range("C41").Value = 'Maven 3rd party JAR'
$ws.Range("D41").Value = 'Guide to installing 3rd party JARs'
$ws.Range("E41").Value = 'Teach'
$ws.Range("F41").Value = 'Succinct'
$ws.Range("G41").Value = 'META-INF, pom file, group-id, artifact-id, packaging, location of local repo, include as dependency'

# ------------------------------------------------------------------
# 3. Add the 'Link' hyperlinks for B26:B41, in row order, so the
#    generated relationship ids line up as rId25..rId40
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B26"), 'http://www.antlr.org/wiki/display/ANTLR4/Getting+Started+with+ANTLR+v4') | Out-Null
$ws.Range("B26").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B27"), 'https://theantlrguy.atlassian.net/wiki/display/ANTLR4/Sample+Maven+plugin+use') | Out-Null
$ws.Range("B27").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B28"), 'https://theantlrguy.atlassian.net/wiki/display/ANTLR4/Getting+Started+with+ANTLR+v4') | Out-Null
$ws.Range("B28").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B29"), 'https://www.jetbrains.com/help/idea/configuring-module-dependencies-and-libraries.html') | Out-Null
$ws.Range("B29").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B30"), 'https://theantlrguy.atlassian.net/wiki/display/ANTLR4/Grammar+Lexicon') | Out-Null
$ws.Range("B30").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B31"), 'https://theantlrguy.atlassian.net/wiki/display/ANTLR4/Grammar+Structure') | Out-Null
$ws.Range("B31").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B32"), 'https://www.digitalocean.com/community/tutorials/how-to-create-an-ssl-certificate-on-nginx-for-ubuntu-14-04') | Out-Null
$ws.Range("B32").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B33"), 'http://stackoverflow.com/questions/3732109/simple-http-server-in-java-using-only-java-se-api') | Out-Null
$ws.Range("B33").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B34"), 'http://www.mkyong.com/java/json-simple-example-read-and-write-json/') | Out-Null
$ws.Range("B34").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B35"), 'http://maven.apache.org/guides/introduction/introduction-to-the-pom.html') | Out-Null
$ws.Range("B35").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B36"), 'http://www.tutorialspoint.com/maven/') | Out-Null
$ws.Range("B36").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B37"), 'http://www.codejava.net/frameworks/maven/how-to-build-project-with-maven') | Out-Null
$ws.Range("B37").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B38"), 'http://www.codejava.net/frameworks/maven/how-to-create-a-java-project-with-maven') | Out-Null
$ws.Range("B38").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B39"), 'http://www.sureshpw.com/2012/09/3-ways-to-run-java-main-from-maven.html') | Out-Null
$ws.Range("B39").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B40"), 'http://stackoverflow.com/questions/4955635/how-to-add-local-jar-files-in-maven-project', 'answer-4955695') | Out-Null
$ws.Range("B40").Value = 'Link'
$ws.Hyperlinks.Add($ws.Range("B41"), 'https://maven.apache.org/guides/mini/guide-3rd-party-jars-local.html') | Out-Null
$ws.Range("B41").Value = 'Link'

# ------------------------------------------------------------------
# 4. Rows 26-32 reuse the same hyperlink text style as rows 2-25 (s=2);
#    copy that formatting across so the style isn't duplicated.
# ------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B26:B32").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 5. Update sheet dimension / selection to match final state
# ------------------------------------------------------------------
$ws.Application.Goto($ws.Range("G40"))
